# "Generate Report for Handoff"
#
# The localization status report moved from "In Translation" to
# "Ready for handoff" and the two "Latest Handoff Datetime" timestamps
# (one per target-locale sheet) were refreshed. Widening the new, longer
# status text also nudges the "Status" / locale-status columns wider on
# every sheet that shows it.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: zh-cn / de-de status columns (E2, F2) -----------------
# G2 ("Latest HO Xliff Generate Date") happens to share its value with the
# de-de sheet's "Latest Handoff Datetime" (H2), so it is refreshed too.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-23 14:46:40"
$wsOverview.Columns.Item(5).ColumnWidth = 16.3
$wsOverview.Columns.Item(6).ColumnWidth = 16.3

# --- zh-cn sheet: Status (C2) + Latest Handoff Datetime (H2) ---------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = "Ready for handoff"
$wsZh.Range("H2").Value = "2016-08-23 14:46:35"
$wsZh.Columns.Item(3).ColumnWidth = 16.3

# --- de-de sheet: Status (C2) + Latest Handoff Datetime (H2) ---------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = "Ready for handoff"
$wsDe.Range("H2").Value = "2016-08-23 14:46:40"
$wsDe.Columns.Item(3).ColumnWidth = 16.3
